$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "62.982.54"
$ws.Range("E2").Value = "  -0.92%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.553.53"
$ws.Range("E3").Value = "  +3.17%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.02%  "

# Row 5 - BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "567.70"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.46%  "

# Row 6 - Solana
$ws.Range("D6").Value = "148.03"
$ws.Range("E6").Value = "  +3.27%  "

# Row 8 - XRP
$ws.Range("E8").Value = "  -2.17%  "

# Row 9 - LidoStakedEther
$ws.Range("D9").Value = "2.552.47"
$ws.Range("E9").Value = "  +3.30%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  -1.99%  "

# Row 11 - Toncoin
$ws.Range("E11").Value = "  -3.13%  "

# Row 12 - TRON
$ws.Range("E12").Value = "  +0.39%  "

# Row 13 - Cardano
$ws.Range("E13").Value = "  -0.70%  "

# Row 14 - Avalanche
$ws.Range("D14").Value = "27.07"

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "3.009.80"
$ws.Range("E15").Value = "  +3.12%  "

# Row 16 - WrappedBTC
$ws.Range("D16").Value = "62.913.36"
$ws.Range("E16").Value = "  -0.82%  "

# Row 17 - ShibaInu
$ws.Range("E17").Value = "  -2.21%  "

# Row 18 - WrappedEther
$ws.Range("D18").Value = "2.519.06"
$ws.Range("E18").Value = "  +1.82%  "

# Row 19 - Chainlink
$ws.Range("D19").Value = "11.49"
$ws.Range("E19").Value = "  +1.33%  "

# Row 20 - BitcoinCash
$ws.Range("D20").Value = "334.42"
$ws.Range("E20").Value = "  -3.06%  "

# Row 21 - Polkadot
$ws.Range("D21").Value = "4.28"
$ws.Range("E21").Value = "  -1.11%  "

# Row 22 - Uniswap
$ws.Range("D22").Value = "6.79"
$ws.Range("E22").Value = "  -0.93%  "

# Row 23 - Dai
$ws.Range("E23").Value = "  +0.05%  "

# Row 24 - Litecoin
$ws.Range("D24").Value = "64.75"
$ws.Range("E24").Value = "  -1.67%  "

# Row 25 - Kaspa
$ws.Range("E25").Value = "  -3.88%  "

# Row 26 - Fetch.AI
$ws.Range("E26").Value = "  +4.10%  "

# Row 27 - Binance-PegBSC-USD
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.09%  "

# Row 28 - SuiNetwork
$ws.Range("E28").Value = "  +11.96%  "

# Row 29 - InternetComputer(DFINITY)
$ws.Range("D29").Value = "8.37"
$ws.Range("E29").Value = "  +1.46%  "

# Row 30 - Aptos
$ws.Range("D30").Value = "7.22"
$ws.Range("E30").Value = "  +4.80%  "

# Row 31 - PEPE
$ws.Range("E31").Value = "  -1.10%  "

# Row 32 - PancakeSwap
$ws.Range("E32").Value = "  +0.08%  "

# Row 33 - Monero
$ws.Range("D33").Value = "177.07"
$ws.Range("E33").Value = "  +0.94%  "

# Row 34 - ImmutableX
$ws.Range("D34").Value = "1.58"
$ws.Range("E34").Value = "  +4.06%  "

# Row 35 - Bittensor
$ws.Range("D35").Value = "408.33"
$ws.Range("E35").Value = "  +9.80%  "

# Row 36 - PolygonEcosystemToken
$ws.Range("D36").Value = "0.398"
$ws.Range("E36").Value = "  -0.81%  "

# Row 37 - EthereumClassic
$ws.Range("D37").Value = "18.83"
$ws.Range("E37").Value = "  -0.92%  "

# Row 40 - Stacks
$ws.Range("E40").Value = "  +0.76%  "

# Row 41 - FirstDigitalUSD
$ws.Range("E41").Value = "  -0.06%  "

# Row 42 - OKB
$ws.Range("D42").Value = "39.22"
$ws.Range("E42").Value = "  -2.77%  "

# Row 43 - Aave
$ws.Range("D43").Value = "151.63"
$ws.Range("E43").Value = "  -0.07%  "

# Row 44 - Filecoin
$ws.Range("D44").Value = "3.74"
$ws.Range("E44").Value = "  +0.06%  "

# Row 45 - InjectiveProtocol
$ws.Range("D45").Value = "20.61"
$ws.Range("E45").Value = "  -0.87%  "

# Row 46 - Mantle
$ws.Range("D46").Value = "0.604"
$ws.Range("E46").Value = "  +0.66%  "

# Row 48 - Hedera
$ws.Range("E48").Value = "  -1.60%  "

# Row 49 - VeChain
$ws.Range("E49").Value = "  +3.97%  "

# Row 50 - EnergySwap
$ws.Range("D50").Value = "18.37"
$ws.Range("E50").Value = "  +0.52%  "

# Row 51 - dogwifhat
$ws.Range("E51").Value = "  +1.82%  "
